# edit.ps1 - apply GuitarPixelmanual.docx revision
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Reposition the tiny "Conector reto 3" straight-line connector shape.
#    posOffset 1019908 -> 1019810 EMU  (=80.3pt,   unchanged-looking but exact)
#    posOffset  302455 ->  -92851 EMU  (=-7.311102362204724pt)
# ---------------------------------------------------------------------------
$shp = $d.Shapes.Item(1)
$shp.Left = 80.3
$shp.Top = -7.311102362204724

# ---------------------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark that sits after "tá doidemais".
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3) "Escolha então qual música gostaria de tocar. Você terá 6 opções, ..."
#    -> split into 3 runs, change 6 -> 4, and re-add "_GoBack" bookmark
#       right after "...terá 4" (before " opções, listadas a seguir: ").
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$oldText = "Escolha então qual música gostaria de tocar. Você terá 6 opções, listadas a seguir: "
$found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Escolha entao' sentence"
}
$sentStart = $find.Parent.Start

$r1 = "Escolha então qual músic"
$r2 = "a gostaria de tocar. Você terá 4"
$r3 = " opções, listadas a seguir: "

$splitA = $sentStart + $r1.Length
$splitB = $splitA + $r2.Length

# Update the word-count digit first (6 -> 4) while the sentence is still one run.
$digitRange = $d.Range($sentStart, $sentStart + $oldText.Length)
$digitRange.Find.Execute("terá 6", $true, $false, $false, $false, $false, $true, 1, $false, "terá 4", 2) | Out-Null

# Split off run1 | run2 (no formatting change -> identical rPr on both pieces,
# exactly mirroring the authored diff).
$rngA = $d.Range($splitA, $splitA)
$d.Bookmarks.Add("TmpSplitA", $rngA) | Out-Null
$d.Bookmarks("TmpSplitA").Delete()

# Insert the "_GoBack" bookmark between run2 and run3 (this also forces the
# run2 | run3 split).
$rngB = $d.Range($splitB, $splitB)
$d.Bookmarks.Add("_GoBack", $rngB) | Out-Null

# ---------------------------------------------------------------------------
# 4) Song list: drop "Carry On My Wayward Son, de Kansas;" and
#    "Bohemian Rhapsody, de Queen." list items; keep "CliFFs OF Dover" item,
#    whose trailing run becomes ", de Eric Johnson." (now the last item).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Carry On My Wayward Son*") {
        $p.Range.Delete()
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Bohemian Rhapsody*") {
        $p.Range.Delete()
        break
    }
}
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute(", de Eric Johnson;", $true, $false, $false, $false, $false, $true, 1, $false, ", de Eric Johnson.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Insert a new empty (center-justified) paragraph right after the
#    paragraph holding the "GUITARRA MANUAL.png" (QWERT keys) picture.
# ---------------------------------------------------------------------------

# Locate the paragraph that contains the "GUITARRA MANUAL.png" drawing: the
# one with an inline picture whose very next paragraph is the "Pressione as
# respectivas teclas..." sentence.
$picParaIndex = $null
$total = $d.Paragraphs.Count
for ($i = 1; $i -le $total; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $nxt = $d.Paragraphs.Item($i + 1)
        if ($nxt.Range.Text -like "*Pressione as respectivas teclas*") {
            $picParaIndex = $i
            break
        }
    }
}
if ($picParaIndex -eq $null) {
    throw "Could not locate the GUITARRA MANUAL picture paragraph"
}

$insertAfterEnd = $d.Paragraphs.Item($picParaIndex).Range.End
$newParaRange = $d.Range($insertAfterEnd, $insertAfterEnd)
$newParaRange.InsertParagraphAfter()

# Format the freshly inserted (now-empty) paragraph: ind firstLine=708, jc=center.
$insertedPara = $d.Paragraphs.Item($picParaIndex + 1)
$insertedPara.Range.ParagraphFormat.FirstLineIndent = 708 / 20.0
$insertedPara.Range.ParagraphFormat.Alignment = 1
